# Auto-generated edit script: applies refreshed market-data values
# to the FFXIV Leve profit-tracking workbook (per commit 'chore: update Sheets via scheduled runner').
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 29000
$ws.Range("J68").Value = 29000
$ws.Range("L68").Value = 29000
$ws.Range("N68").Value = -30498

$ws.Range("H71").Value = 29000
$ws.Range("J71").Value = 29000
$ws.Range("L71").Value = 87000
$ws.Range("N71").Value = -94488

$ws.Range("H137").Value = 5650.5
$ws.Range("I137").Value = 6784
$ws.Range("J137").Value = 3836.9
$ws.Range("K137").Value = 20352
$ws.Range("L137").Value = 11510.7
$ws.Range("M137").Value = -17802
$ws.Range("N137").Value = -16610.7

$ws.Range("H141").Value = 393802.97
$ws.Range("I141").Value = 1314.4615
$ws.Range("J141").Value = 606400.9399999999
$ws.Range("K141").Value = 3943.3845
$ws.Range("L141").Value = 1819202.82
$ws.Range("M141").Value = 1236.6155
$ws.Range("N141").Value = -1829562.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 32611992
$ws.Range("I2").Value = 107145860
$ws.Range("J2").Value = 3425.8125
$ws.Range("K2").Value = 107145860
$ws.Range("L2").Value = 3425.8125
$ws.Range("M2").Value = -107145747
$ws.Range("N2").Value = -3651.8125

$ws.Range("H21").Value = 43793
$ws.Range("I21").Value = 1625
$ws.Range("J21").Value = 100017
$ws.Range("K21").Value = 1625
$ws.Range("L21").Value = 100017
$ws.Range("M21").Value = -1251
$ws.Range("N21").Value = -100765

$ws.Range("H109").Value = 30200
$ws.Range("J109").Value = 30200
$ws.Range("L109").Value = 30200
$ws.Range("N109").Value = -32974

$ws.Range("H116").Value = 32611992
$ws.Range("I116").Value = 107145860
$ws.Range("J116").Value = 3425.8125
$ws.Range("K116").Value = 107145860
$ws.Range("L116").Value = 3425.8125
$ws.Range("M116").Value = -107143566
$ws.Range("N116").Value = -8013.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 32611992
$ws.Range("I3").Value = 107145860
$ws.Range("J3").Value = 3425.8125
$ws.Range("K3").Value = 107145860
$ws.Range("L3").Value = 3425.8125
$ws.Range("M3").Value = -107145746
$ws.Range("N3").Value = -3653.8125

$ws.Range("H134").Value = 4325.775
$ws.Range("I134").Value = 4027.1316
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 12081.3948
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -9546.3948
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1516.8462
$ws.Range("I16").Value = 792.63635
$ws.Range("K16").Value = 792.63635
$ws.Range("M16").Value = -505.63635

$ws.Range("H31").Value = 3147.7302
$ws.Range("I31").Value = 2230.5625
$ws.Range("J31").Value = 3662.6316
$ws.Range("K31").Value = 2230.5625
$ws.Range("L31").Value = 3662.6316
$ws.Range("M31").Value = -1935.5625
$ws.Range("N31").Value = -4252.631600000001

$ws.Range("H34").Value = 3147.7302
$ws.Range("I34").Value = 2230.5625
$ws.Range("J34").Value = 3662.6316
$ws.Range("K34").Value = 2230.5625
$ws.Range("L34").Value = 3662.6316
$ws.Range("M34").Value = -2028.5625
$ws.Range("N34").Value = -4066.6316

$ws.Range("H58").Value = 16671329
$ws.Range("I58").Value = 2760.375
$ws.Range("J58").Value = 35721120
$ws.Range("K58").Value = 2760.375
$ws.Range("L58").Value = 35721120
$ws.Range("M58").Value = -2557.375
$ws.Range("N58").Value = -35721526

$ws.Range("H113").Value = 1516.8462
$ws.Range("I113").Value = 792.63635
$ws.Range("K113").Value = 792.63635
$ws.Range("M113").Value = 1377.36365

$ws.Range("H123").Value = 36268
$ws.Range("J123").Value = 36268
$ws.Range("L123").Value = 36268
$ws.Range("N123").Value = -46068

$ws.Range("H124").Value = 31450
$ws.Range("J124").Value = 31450
$ws.Range("L124").Value = 31450
$ws.Range("N124").Value = -36360

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 2222.4866
$ws.Range("I132").Value = 1874.8125
$ws.Range("J132").Value = 4447.6
$ws.Range("K132").Value = 5624.4375
$ws.Range("L132").Value = 13342.8
$ws.Range("M132").Value = -3094.4375
$ws.Range("N132").Value = -18402.8

$ws.Range("H134").Value = 12197401
$ws.Range("I134").Value = 12822447
$ws.Range("K134").Value = 38467341
$ws.Range("M134").Value = -38464806

$ws.Range("H136").Value = 16671329
$ws.Range("I136").Value = 2760.375
$ws.Range("J136").Value = 35721120
$ws.Range("K136").Value = 8281.125
$ws.Range("L136").Value = 107163360
$ws.Range("M136").Value = -5731.125
$ws.Range("N136").Value = -107168460

$ws.Range("H141").Value = 25105.883
$ws.Range("J141").Value = 25105.883
$ws.Range("L141").Value = 25105.883
$ws.Range("N141").Value = -35465.883

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 9921.764999999999
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 10473.125
$ws.Range("K22").Value = 3300
$ws.Range("L22").Value = 31419.375
$ws.Range("M22").Value = -3131
$ws.Range("N22").Value = -31757.375

$ws.Range("H27").Value = 9921.764999999999
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 10473.125
$ws.Range("K27").Value = 3300
$ws.Range("L27").Value = 31419.375
$ws.Range("M27").Value = -3198
$ws.Range("N27").Value = -31623.375

$ws.Range("H36").Value = 2400.4
$ws.Range("I36").Value = 2
$ws.Range("K36").Value = 6
$ws.Range("M36").Value = 163

$ws.Range("H111").Value = 2566
$ws.Range("I111").Value = 415
$ws.Range("K111").Value = 1245
$ws.Range("M111").Value = 1822

$ws.Range("H113").Value = 839.1
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 839.1
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2517.3
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6857.3

$ws.Range("H131").Value = 1223.0814
$ws.Range("I131").Value = 2032.3077
$ws.Range("J131").Value = 1078.9727
$ws.Range("K131").Value = 6096.9231
$ws.Range("L131").Value = 3236.9181
$ws.Range("M131").Value = -1056.9231
$ws.Range("N131").Value = -13316.9181

$ws.Range("H134").Value = 2248.625
$ws.Range("I134").Value = 1664.8334
$ws.Range("K134").Value = 4994.5002
$ws.Range("M134").Value = 75.4997999999996

$ws.Range("H137").Value = 3139.9412
$ws.Range("J137").Value = 5000
$ws.Range("L137").Value = 15000
$ws.Range("N137").Value = -25200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4057.9
$ws.Range("I80").Value = 3614.3635
$ws.Range("J80").Value = 4600
$ws.Range("K80").Value = 3614.3635
$ws.Range("L80").Value = 4600
$ws.Range("M80").Value = -2616.3635
$ws.Range("N80").Value = -6596

$ws.Range("H83").Value = 4057.9
$ws.Range("I83").Value = 3614.3635
$ws.Range("J83").Value = 4600
$ws.Range("K83").Value = 18071.8175
$ws.Range("L83").Value = 23000
$ws.Range("M83").Value = -13079.8175
$ws.Range("N83").Value = -32984

$ws.Range("H132").Value = 3544.675
$ws.Range("I132").Value = 4440.4707
$ws.Range("J132").Value = 2882.5652
$ws.Range("K132").Value = 13321.4121
$ws.Range("L132").Value = 8647.695599999999
$ws.Range("M132").Value = -10791.4121
$ws.Range("N132").Value = -13707.6956

$ws.Range("H133").Value = 29000
$ws.Range("J133").Value = 29000
$ws.Range("L133").Value = 29000
$ws.Range("N133").Value = -39120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4711.069
$ws.Range("I136").Value = 4226
$ws.Range("K136").Value = 12678
$ws.Range("M136").Value = -10128

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8594.166999999999
$ws.Range("J41").Value = 8594.166999999999
$ws.Range("L41").Value = 8594.166999999999
$ws.Range("N41").Value = -9374.166999999999

$ws.Range("H58").Value = 14000
$ws.Range("J58").Value = 14000
$ws.Range("L58").Value = 14000
$ws.Range("N58").Value = -14616

$ws.Range("H98").Value = 29999.666
$ws.Range("J98").Value = 29999.666
$ws.Range("L98").Value = 29999.666
$ws.Range("N98").Value = -35989.666
